$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "modified" timestamp in B20 (new .ttl regenerated from the
# Google sheet later the same day).
$ws.Range("B20").Value = "2022-06-20T10:51:30+00:00"

# New vocabulary terms appended to the sheet (rows 25-30), mirroring the
# existing "Identifier / prefLabel / .../ broader" layout (columns A, B, F).
# Entry format: @(identifier, prefLabel, broaderTermLabel)
$newRows = @(
    @("vocab:1002", "reaction time at lexical decision task (without prospective memory component)", "reaction time"),
    @("vocab:1003", "reaction time at lexical decision task (with prospective memory component)", "reaction time"),
    @("vocab:1004", "accuracy", ""),
    @("vocab:1005", "accuracy at prospective memory task", "accuracy"),
    @("vocab:1006", "accuracy at lexical decision task (without prospective memory component)", "accuracy"),
    @("vocab:1007", "accuracy at lexical decision task (with prospective memory component)", "accuracy")
)

$startRow = 25
$lastCol = 37  # Column AK, matching the sheet's existing A:AK extent

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Touch every column A..AK so the row keeps the sheet's established
    # rectangular shape (every other row has all 37 columns materialised,
    # even where blank).
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Style = "Normal"
    }

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    if ($rowData[2] -ne "") {
        $ws.Cells.Item($r, 6).Value = $rowData[2]
    }
}
